$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '[''XS1819680288'', ''XS1819680528'', ''XS2083302419'', ''XS2083302500'', ''XS2446175577'']'
$ws.Range("D5").Value = '[''XS1558078496'', ''XS1775617464'', ''XS1807305328'', ''XS1953057061'', ''XS1953057491'', ''XS1980065301'', ''XS1980255936'', ''XS2079842642'', ''XS2176899701'', ''XS2241075014'', ''XS2297220423'', ''XS2297221405'', ''XS2297226545'', ''XS2391395154'', ''XS2391398174'', ''XS2530049837'']'
$ws.Range("D16").Value = '[''XS1405777589'', ''XS1575967218'', ''XS1575968026'', ''XS1750113661'', ''XS1750114396'', ''XS1944412664'', ''XS1944412748'', ''XS2234859101'', ''XS2248458395'', ''XS2288906857'', ''XS2310799809'', ''XS2351310482'', ''XS2689095086'']'
$ws.Range("E2").Value = '[''8,25 ANGOL 09-05-2028 (REGS)'', ''REPUBLIC OF ANGOLA 8.25% 09.05.2028'', ''REPUBLIC OF ANGOLA 9.375% 08.05.2048'', ''REPUBLIC OF ANGOLA 8% 26.11.2029'', ''9,125 ANGOL 26-11-2049 (REGS)'', ''REPUBLIC OF ANGOLA 9.125% 26.11.2049'', ''8,75 ANGOL 14-04-2032 (REGS)'', ''REPUBLIC OF ANGOLA 8.75% 14.04.2032'']'
$ws.Range("E5").Value = '[''ARAB REPUBLIC OF EGYPT 8.5% 31.01.2047'', ''7.90% Egypt Government International Bond 2048'', ''ARAB REPUBLIC OF EGYPT 7.903% 21.02.2048'', ''ARAB REPUBLIC OF EGYPT 5.625% 16.04.2030'', ''ARAB REPUBLIC OF EGYPT 7.6003% 01.03.2029'', ''8.70% Egypt Government International Bond 2049'', ''EGYPT 8.7002 03/01/49'', ''ARAB REPUBLIC OF EGYPT 4.75% 11.04.2025'', ''ARAB REPUBLIC OF EGYPT 6.375% 11.04.2031'', ''EGYPT 6 3/8 04/11/31'', ''ARAB REPUBLIC OF EGYPT 7.0529% 15.01.2032'', ''8,875 EGYPT 29-05-2050 (REGS)'', ''ARAB REPUBLIC OF EGYPT 8.875% 29.05.2050'', ''Egypt, Arab Republic Of (Government) 5.25% 20251006'', ''ARAB REPUBLIC OF EGYPT 3.875% 16.02.2026'', ''ARAB REPUBLIC OF EGYPT 7.5% 16.02.2061'', ''5,875 EGYPT 16-02-2031 (REGS)'', ''5.88% Egypt Government International Bond 2031'', ''7.30% Egypt Government International Bond 2033'', ''8.75% Egypt Government International Bond 2051'', ''EGYPT 8 3/4 09/30/51'', ''EGYPT TASKEEK CO 10.875% 28.02.2026'', ''EGYPT TASKEEK COMPANY 10.875% 28.02.2026'']'
$ws.Range("E16").Value = '[''4,75 OMAN 15-06-2026 (REGS)'', ''5.38% Oman Government International Bond 2027'', ''6,5 OMAN 08-03-2047 (REGS)'', ''OMAN GOV INTERNTL BOND 6.5% 08.03.2047'', ''5.63% Oman Government International Bond 2028'', ''6,75 OMAN 17-01-2048 (REGS)'', ''OMAN 6 3/4 01/17/48'', ''OMAN GOV INTERNTL BOND 6.75% 17.01.2048'', ''OMAN 4 7/8 02/01/25'', ''6 OMAN 01-08-2029 (REGS)'', ''6.00% Oman Government International Bond 2029'', ''6.75% Oman Government International Bond 2027'', ''OQ SAOC 5.125% 06.05.2028'', ''OMAN GOV INTERNTL BOND 7% 25.01.2051'', ''BANK MUSCAT SAOG 4.75% 17.03.2026'', ''OMAN SOVEREIGN SUKUK 4.875% 15.06.2030'', ''EDO SUKUK LTD 5.875% 21.09.2033'']'
$ws.Range("E18").Value = '[''3,25 QATAR 02-06-2026 (REGS)'', ''4,5 QATAR 23-04-2028 (REGS)'', ''5.10% Qatar Government International Bond 2048'', ''STATE OF QATAR 5.103% 23.04.2048'', ''STATE OF QATAR 4.817% 14.03.2049'', ''STATE OF QATAR 3.4% 16.04.2025'', ''3,75 QATAR 16-04-2030 (REGS)'', ''4,4 QATAR 16-04-2050 (REGS)'', ''STATE OF QATAR 4.4% 16.04.2050'', ''QATAR ENERGY 2.25% 12.07.2031'']'
$ws.Range("F2").Value = '[''Republic of Angola'', ''Angola'', ''Angolan Government International Bond'', ''REPUBLIC OF ANGOLA  '']'
$ws.Range("F16").Value = '[''Sultanate of Oman'', ''5.38% Oman Government International Bond 2027'', ''Oman Government International Bond'', ''Oman'', ''OMAN GOV INTERNTL BOND  '', ''5.63% Oman Government International Bond 2028'', ''Oman Government International Bond 17.01.2048'', ''Oman Government International Obligation 17.01.2048'', ''Oman Government International Bond 01.02.2025'', ''Oman Government International Obligation 01.02.2025'', ''6.00% Oman Government International Bond 2029'', ''6.75% Oman Government International Bond 2027'', ''OMAN SOVEREIGN SUKUK  '', ''Oman Sovereign Sukuk Co'']'
$ws.Range("F18").Value = '[''State of Qatar'', ''5.10% Qatar Government International Bond 2048'', ''Qatar Government International Bond'', ''Qatar'', ''STATE OF QATAR  '', ''QATAR ENERGY  '', ''QatarEnergy'']'
$ws.Range("G2").Value = '[''Herning'', ''Fanø'', ''Kalundborg'', ''Rødovre'', ''Vejen'']'
$ws.Range("G16").Value = '[''Herning'', ''Ringkøbing-Skjern'', ''Region Nordjylland'', ''Kalundborg'', ''Rødovre'', ''Vejen'', ''Greve'', ''Kolding'', ''Skive'', ''Fredericia'', ''Hillerød'', ''Lemvig'', ''Varde'', ''Aabenraa'', ''Fanø'']'
